$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Summary block (rows 10-12): fill in the actual grading numbers instead of
# the all-zero / "Absent" placeholders, and give the row-label cells in
# column A the same boxed "mtitleStyle" used by the header row above them.
# ---------------------------------------------------------------------------
$ws.Range("A10").Style = "mtitleStyle"
$ws.Range("B10").Value = 12
$ws.Range("C10").Value = 3
$ws.Range("D10").Value = 13
$ws.Range("E10").Value = 28

$ws.Range("A11").Style = "mtitleStyle"
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1
$ws.Range("D11").Value = 0

$ws.Range("A12").Style = "mtitleStyle"
$ws.Range("B12").Value = 48
$ws.Range("C12").Value = -3
$ws.Range("E12").Value = "45/112"

# ---------------------------------------------------------------------------
# The sheet used to show three side-by-side Student/Correct Ans blocks
# (A:B, D:E, G:H). The third block (columns G:H) is no longer needed, so
# remove those columns outright.
# ---------------------------------------------------------------------------
$ws.Range("G:H").EntireColumn.Delete()

# ---------------------------------------------------------------------------
# First Student/Correct Ans block (column A): fill in the student's actual
# answers, colour-coded with the workbook's own named styles -
# correctStyle (matches the correct answer), incorrectStyle (doesn't match)
# - leaving the "not attempted" rows with their existing blank normalStyle.
# ---------------------------------------------------------------------------
$ws.Range("A16").Style = "correctStyle"
$ws.Range("A16").Value = "Option A"

$ws.Range("A18").Style = "correctStyle"
$ws.Range("A18").Value = "Option B"

$ws.Range("A19").Style = "correctStyle"
$ws.Range("A19").Value = "Option C"

$ws.Range("A21").Style = "correctStyle"
$ws.Range("A21").Value = "Option C"

$ws.Range("A22").Style = "incorrectStyle"
$ws.Range("A22").Value = "Option A"

$ws.Range("A23").Style = "correctStyle"
$ws.Range("A23").Value = "Option D"

$ws.Range("A25").Style = "correctStyle"
$ws.Range("A25").Value = "Option A"

$ws.Range("A26").Style = "correctStyle"
$ws.Range("A26").Value = "Option C"

$ws.Range("A27").Style = "incorrectStyle"
$ws.Range("A27").Value = "Option D"

$ws.Range("A28").Style = "correctStyle"
$ws.Range("A28").Value = "Option D"

$ws.Range("A32").Style = "correctStyle"
$ws.Range("A32").Value = "Option C"

$ws.Range("A35").Style = "incorrectStyle"
$ws.Range("A35").Value = "Option A"

$ws.Range("A39").Style = "correctStyle"
$ws.Range("A39").Value = "Option D"

# ---------------------------------------------------------------------------
# Second Student/Correct Ans block (columns D:E): only rows 16-18 still
# apply to this quiz, so fill those two student-answer cells in and clear
# the remaining rows (19-40) completely - they no longer have a question.
# ---------------------------------------------------------------------------
$ws.Range("D16").Style = "correctStyle"
$ws.Range("D16").Value = "Option A"

$ws.Range("D18").Style = "correctStyle"
$ws.Range("D18").Value = "Option D"

$ws.Range("D19:E40").Clear()
